# "Draft elec sector edits"
#
# The "BPaFF-BITPTaP" sheet (Is-This-Plant-Type-a-Peaker flags) had its
# natural-gas-nonpeaker / nuclear / hydro flags (B2:B4) flipped from 1 to 0,
# matching the all-zero pattern already used on the "BPaFF-BDTPTPF"
# (Does-This-Plant-Type-Provide-Flexibility) sheet. The dependent formula
# cells (B13:B17) recalculate automatically since they just reference B2,
# B6, B9 and B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPaFF-BITPTaP")

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0

# Restore "About" as the active/front sheet (it was the selected tab in the
# edited workbook).
$about = $wb.Worksheets.Item("About")
$about.Activate()
$about.Range("A1").Select()
